# Insert a new survey-metadata row for EB 95.1 (March-April 2021) above the
# existing most-recent entry (ZA7750 / 94.2), shifting all following rows
# down by one, then move the active selection to D3 (matching the author's
# final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 2 (and everything below it) down by one row.
$ws.Rows(2).Insert() | Out-Null

# Populate the newly-opened row 2 with the new survey entry.
$ws.Range("A2").Value2 = "ZA7781"
$ws.Range("B2").Value2 = "'95.1"
$ws.Range("C2").Value2 = "March-April 2021"
$ws.Range("D2").Value2 = "European Parliament Spring Survey, Climate Change, and EU Consumer Habits Regarding Fishery and Aquaculture Products"

# Match the saved selection state from the edit.
$ws.Range("D3").Select() | Out-Null
